# Code_Diagram.pptx slide 1 edit:
#  - shift several boxes/connectors horizontally (and a couple vertically)
#  - rewrite the "summary_stats_continuous_and_likert" box's description text
#  - remove the now-unused connector between shapes 4 and 28 (id 29)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$emuPerPt = 12700
# The host stores Left/Top/Width/Height as single-precision points and
# truncates (rather than rounds) when converting back to EMU on save, so a
# plain EMU/12700 division can land one EMU below the intended value. Adding
# half an EMU (in points) before assignment compensates for that truncation
# and reproduces the exact target EMU values.
$emuRoundingFudge = 0.5 / $emuPerPt

function Emu-ToPt($emu) {
    return ($emu / $emuPerPt) + $emuRoundingFudge
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# Reposition shapes (offsets taken straight from the OOXML <a:off> EMU values,
# converted to points as PowerPoint's Shape.Left/Top expect).
$moves = @(
    @{ Id = 4;  X = 809175;  Y = 2860555 },
    @{ Id = 5;  X = 213687;  Y = 204124 },
    @{ Id = 6;  X = 3790233; Y = 204124 },
    @{ Id = 7;  X = 3790233; Y = 1867593 },
    @{ Id = 8;  X = 3790233; Y = 3531062 },
    @{ Id = 10; X = 3035139; Y = 933797 },
    @{ Id = 12; X = 3035139; Y = 2597266 },
    @{ Id = 15; X = 3035139; Y = 3590228 },
    @{ Id = 23; X = 3790233; Y = 5194531 },
    @{ Id = 24; X = 3035139; Y = 3590228 },
    @{ Id = 27; X = 6771291; Y = 204123 },
    @{ Id = 28; X = 213687;  Y = 5194531 },
    @{ Id = 32; X = 6016197; Y = 933797 },
    @{ Id = 41; X = 6016197; Y = 2597266 },
    @{ Id = 43; X = 6016197; Y = 4260735 },
    @{ Id = 46; X = 6016197; Y = 5924203 },
    @{ Id = 48; X = 9752349; Y = 2699328 },
    @{ Id = 49; X = 8997255; Y = 3429000 }
)

foreach ($mv in $moves) {
    $sh = Get-ShapeById $shapes $mv.Id
    $sh.Left = Emu-ToPt $mv.X
    $sh.Top  = Emu-ToPt $mv.Y
}

# Update the description text in the "summary_stats_continuous_and_likert" box (id 28).
$summaryShape = Get-ShapeById $shapes 28
$descParagraph = $summaryShape.TextFrame.TextRange.Paragraphs(2)
$descParagraph.Runs(1).Text = "Creates a fancy table with the summary of each attribute of the entire  training set."

# Remove the now-obsolete connector from shape 4 to shape 28 (id 29).
$oldConnector = Get-ShapeById $shapes 29
if ($oldConnector -ne $null) {
    $oldConnector.Delete()
}
